$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.541.23"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.45%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.856.00"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.14%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.0000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.32%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'233.82"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.42%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.27%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4748"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +3.39%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2748"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +3.27%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06322"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.84%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'17.75"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +11.97%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = "'TRON"
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').Value = "'0.07453"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.59%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = "'WrappedEther"
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = "'1.798.11"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -2.03%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.988"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +3.03%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'84.69"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +2.34%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.6275"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.28%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'30.505.39"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.62%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'246.46"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +9.79%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.18%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'12.71"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +3.17%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.000007328"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +2.05%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.9995"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.66%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'4.941"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +2.14%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'5.926"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.78%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'9.127"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.83%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'162.70"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.04%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'18.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.37%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'1.878"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.18%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'0.1023"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.18%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'1.350"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -2.08%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'4.027"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.09%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'3.835"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +2.57%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.04844"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.08%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'1.137"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.11%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.7022"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.06%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.699"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.39%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.01900"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +5.23%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +3.37%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'2.006"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +5.07%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.8760"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.91%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'106.82"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +3.74%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.28%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'5.548"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.22%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.4057"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +2.26%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'7.213"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +5.36%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'62.85"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +6.94%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.1205"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +3.03%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'33.67"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +4.25%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'8.533"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +2.11%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.31%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.352"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.04%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.3690"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.59%  "
$ws.Range('E51').Style = 'Normal'
